$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update H2:H11 values from 60 to 59
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 59
}

# Update the active selection from H19 to H17
$ws.Range("H17").Select()
